$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Transmittals_New")
for ($row = 3; $row -le 19; $row++) {
    $cell = $ws1.Range("M$row")
    $cell.Formula = "=CONCATENATE(ROW()-1," + [char]34 + " of " + [char]34 + ",COUNTA(A2:A100)," + [char]34 + " " + [char]34 + ",N$row)"
}

$ws2 = $wb.Worksheets.Item("Transmittals_New_ActionRequired")
$cell2 = $ws2.Range("M3")
$cell2.Formula = "=CONCATENATE(ROW()-1," + [char]34 + " of " + [char]34 + ",COUNTA(A2:A100)," + [char]34 + " " + [char]34 + ",N3)"

$excel.Calculate()
